$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly data refresh shifts every existing "Berenjena" price record
# (rows 22-38) down by one row, and inserts a brand-new record at row 22
# with the latest date. Inserting a whole row at position 22 performs
# exactly that shift (and carries the row's formatting, e.g. the date
# number format on column D) for rows 22-38 -> 23-39 in one step.
$ws.Rows("22:22").Insert()

# Populate the newly inserted row 22 with the latest week's record.
$ws.Range("A22").Value = 11
$ws.Range("B22").Value = "Vega Monumental Concepción"
$ws.Range("C22").Value = "Bíobío"
$ws.Range("D22").Value2 = 44467
$ws.Range("E22").Value = 8
$ws.Range("F22").Value = 100112001
$ws.Range("G22").Value = "Berenjena"
$ws.Range("H22").Value = "Sin especificar"
$ws.Range("I22").Value = "Primera"
$ws.Range("J22").Value = 100
$ws.Range("K22").Value = 9000
$ws.Range("L22").Value = 10000
$ws.Range("M22").Value = 9500
$ws.Range("N22").Value = "$/caja 60 unidades"
$ws.Range("O22").Value = "Región de Arica y Parinacota"
$ws.Range("P22").Value = 158
$ws.Range("Q22").Value = 60
$ws.Range("R22").Value = "Hortaliza"
